$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiments")
$ws.Activate()

# Insert 6 new rows above row 43 (5 new data rows + 1 blank spacer row),
# so old rows 43:46 become 49:52
$ws.Range("A43:A48").EntireRow.Insert()

# New row 43: DORN (Weighted histogram matching) - Intensity Only
$ws.Range("A43").Value = "DORN (Weighted histogram matching)"
$ws.Range("B43").Value = 0.90449999999999997
$ws.Range("C43").Value = 0.97050000000000003
$ws.Range("D43").Value = 0.98919999999999997
$ws.Range("F43").Value = 0.41420000000000001
$ws.Range("G43").Value = 0.0912
$ws.Range("I43").Value = 0.0395
$ws.Range("L43").Value = "Intensity Only"

# New row 44: DORN (Weighted histogram matching) - Intensity and Falloff
$ws.Range("A44").Value = "DORN (Weighted histogram matching)"
$ws.Range("B44").Value = 0.90449999999999997
$ws.Range("C44").Value = 0.97050000000000003
$ws.Range("D44").Value = 0.98919999999999997
$ws.Range("F44").Value = 0.4143
$ws.Range("G44").Value = 0.0912
$ws.Range("I44").Value = 0.0395
$ws.Range("L44").Value = "Intensity and Falloff"

# New row 45: DORN (Weighted histogram matching) - Intensity, Falloff, and DC/Ambient
$ws.Range("A45").Value = "DORN (Weighted histogram matching)"
$ws.Range("B45").Value = 0.90410000000000001
$ws.Range("C45").Value = 0.97
$ws.Range("D45").Value = 0.9889
$ws.Range("F45").Value = 0.4173
$ws.Range("G45").Value = 0.0902
$ws.Range("I45").Value = 0.0396
$ws.Range("L45").Value = "Intensity, Falloff, and DC/Ambient"

# New row 46: DORN (Weighted histogram matching) - Intensity, Falloff, DC/Ambient, and Jitter
$ws.Range("A46").Value = "DORN (Weighted histogram matching)"
$ws.Range("B46").Value = 0.90410000000000001
$ws.Range("C46").Value = 0.97009999999999996
$ws.Range("D46").Value = 0.9889
$ws.Range("F46").Value = 0.4168
$ws.Range("G46").Value = 0.0903
$ws.Range("I46").Value = 0.0396
$ws.Range("L46").Value = "Intensity, Falloff, DC/Ambient, and Jitter"

# New row 47: DORN (Weighted histogram matching) - Intensity, Falloff, DC/Ambient, Jitter, and Poisson Noise
$ws.Range("A47").Value = "DORN (Weighted histogram matching)"
$ws.Range("B47").Value = 0.90310000000000001
$ws.Range("C47").Value = 0.96909999999999996
$ws.Range("D47").Value = 0.98809999999999998
$ws.Range("F47").Value = 0.45860000000000001
$ws.Range("G47").Value = 0.0926
$ws.Range("I47").Value = 0.04
$ws.Range("L47").Value = "Intensity, Falloff, DC/Ambient, Jitter, and Poisson Noise"

# Row heights for the newly inserted data rows (17px, matching the single-line rows elsewhere)
$ws.Rows("43:47").RowHeight = 17

# Row 48 is left blank as a visual spacer (matches original gap before the summary rows)

# Update view/selection to match final state: select whole row 52, scroll so row 27 is at top
$ws.Rows(52).Select()
$excel.ActiveWindow.ScrollRow = 27
$excel.ActiveWindow.ScrollColumn = 1
